$wb = $excel.ActiveWorkbook

# Update both the "展览" and "全部类型" sheets which contain the same data table.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 493
    $ws.Range("F5").Value = 4880
}
